$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Adjust continuous skills' period values (row 14 and row 17), and the
# delay value on row 17, per the commit: "Adjust some continuous skills' period"
$ws.Range("E14").Value = 2
$ws.Range("E17").Value = 2
$ws.Range("H17").Value = 2

# Update the active selection to match the saved view state
$ws.Range("I25").Select()
